$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 190.2720336914062
$ws.Range("B3").Value = 188.2866973876953
$ws.Range("B4").Value = 179.1955718994141
$ws.Range("B5").Value = 179.2254791259766
$ws.Range("B6").Value = 181.0870666503906
$ws.Range("B7").Value = 175.9044952392578
$ws.Range("B8").Value = 178.6998138427734
$ws.Range("B9").Value = 180.3567657470703
$ws.Range("B10").Value = 195.4346923828125
$ws.Range("B11").Value = 206.8852233886719
$ws.Range("B12").Value = 211.4739837646484
$ws.Range("B13").Value = 242.8189544677734
$ws.Range("B14").Value = 263.66455078125
$ws.Range("B15").Value = 260.4147033691406
$ws.Range("B16").Value = 215.7895202636719
$ws.Range("B17").Value = 196.5757293701172
$ws.Range("B18").Value = 177.3250885009766
$ws.Range("B19").Value = 168.6316833496094
$ws.Range("B20").Value = 154.1441345214844
$ws.Range("B21").Value = 142.3265075683594
$ws.Range("B22").Value = 141.2400817871094
$ws.Range("B23").Value = 130.4773254394531
$ws.Range("B24").Value = 123.952522277832
$ws.Range("B25").Value = 118.3817138671875
